# Refresh the cryptocurrency price/volume table on Sheet1.
# Most rows only get updated Price (D) / Volume(1h) (E) figures; a few rows
# (19/20, 42/43/44, 50/51) also swap which coin occupies which row, carrying
# the coin name (B), link (C), price (D) and volume (E) with them.
#
# Several "Price" values in column D are plain text that happens to look like
# a number (e.g. "18.98"). A bare $ws.Range(...).Value = "18.98" assignment
# would get auto-coerced to a numeric cell by Excel, which would not match the
# original text-cell formatting. To keep those cells as text (same as the rest
# of the column), we briefly force the cell to a text NumberFormat, write the
# value, then restore the "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    if ($value -match "^[+-]?[0-9]*\.?[0-9]+$") {
        # Looks like a plain number to Excel -- pin it down as text so the
        # stored cell keeps its original (string) type instead of becoming numeric.
        $rng.NumberFormat = "@"
        $rng.Value = $value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $value
    }
}


# Row 2
Set-TextValue "D2" "26.720.79"
Set-TextValue "E2" "  -0.17%  "

# Row 3
Set-TextValue "D3" "1.633.60"
Set-TextValue "E3" "  -0.86%  "

# Row 4
Set-TextValue "E4" "  +0.00%  "

# Row 5
Set-TextValue "D5" "218.54"
Set-TextValue "E5" "  +0.86%  "

# Row 6
Set-TextValue "E6" "  -1.47%  "

# Row 7
Set-TextValue "E7" "  +0.00%  "

# Row 8
Set-TextValue "E8" "  -1.22%  "

# Row 9
Set-TextValue "E9" "  -1.06%  "

# Row 10
Set-TextValue "D10" "18.98"
Set-TextValue "E10" "  -1.46%  "

# Row 11
Set-TextValue "D11" "0.0842"
Set-TextValue "E11" "  -0.16%  "

# Row 12
Set-TextValue "D12" "1.860.57"
Set-TextValue "E12" "  -0.87%  "

# Row 13
Set-TextValue "D13" "1.650.88"
Set-TextValue "E13" "  +0.20%  "

# Row 14
Set-TextValue "E14" "  -2.43%  "

# Row 15
Set-TextValue "E15" "  -1.93%  "

# Row 16
Set-TextValue "D16" "64.08"
Set-TextValue "E16" "  -2.41%  "

# Row 17
Set-TextValue "D17" "26.701.94"
Set-TextValue "E17" "  -0.29%  "

# Row 18
Set-TextValue "E18" "  -2.85%  "

# Row 19
Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "211.45"
Set-TextValue "E19" "  -2.64%  "

# Row 20
Set-TextValue "B20" "Dai"
Set-TextValue "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.01"
Set-TextValue "E20" "  +0.05%  "

# Row 21
Set-TextValue "E21" "  -1.34%  "

# Row 22
Set-TextValue "D22" "6.19"
Set-TextValue "E22" "  -2.01%  "

# Row 23
Set-TextValue "D23" "2.32"
Set-TextValue "E23" "  -4.22%  "

# Row 24
Set-TextValue "D24" "9.19"
Set-TextValue "E24" "  -2.90%  "

# Row 25
Set-TextValue "D25" "146.97"
Set-TextValue "E25" "  +0.91%  "

# Row 26
Set-TextValue "E26" "  -0.06%  "

# Row 27
Set-TextValue "E27" "  -2.26%  "

# Row 28
Set-TextValue "D28" "7.00"
Set-TextValue "E28" "  -2.99%  "

# Row 29
Set-TextValue "D29" "15.53"
Set-TextValue "E29" "  -1.89%  "

# Row 30
Set-TextValue "E30" "  -4.24%  "

# Row 31
Set-TextValue "E31" "  +0.76%  "

# Row 32
Set-TextValue "D32" "3.36"
Set-TextValue "E32" "  +0.40%  "

# Row 33
Set-TextValue "E33" "  -2.51%  "

# Row 34
Set-TextValue "D34" "1.264.29"
Set-TextValue "E34" "  -0.93%  "

# Row 35
Set-TextValue "E35" "  -1.79%  "

# Row 36
Set-TextValue "E36" "  +0.29%  "

# Row 37
Set-TextValue "D37" "0.0173"
Set-TextValue "E37" "  -3.16%  "

# Row 38
Set-TextValue "D38" "0.524"
Set-TextValue "E38" "  -3.32%  "

# Row 39
Set-TextValue "E39" "  +0.02%  "

# Row 40
Set-TextValue "D40" "0.803"
Set-TextValue "E40" "  -3.67%  "

# Row 41
Set-TextValue "D41" "0.799"
Set-TextValue "E41" "  -2.35%  "

# Row 42
Set-TextValue "B42" "FraxShare"
Set-TextValue "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.26"
Set-TextValue "E42" "  -3.20%  "

# Row 43
Set-TextValue "B43" "MXToken"
Set-TextValue "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.15"
Set-TextValue "E43" "  -4.37%  "

# Row 44
Set-TextValue "B44" "RocketPoolETH"
Set-TextValue "C44" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D44" "1.771.58"
Set-TextValue "E44" "  -1.52%  "

# Row 45
Set-TextValue "D45" "91.45"
Set-TextValue "E45" "  -0.76%  "

# Row 46
Set-TextValue "D46" "59.68"
Set-TextValue "E46" "  +0.14%  "

# Row 47
Set-TextValue "D47" "1.57"
Set-TextValue "E47" "  -3.33%  "

# Row 48
Set-TextValue "D48" "0.0515"
Set-TextValue "E48" "  -0.16%  "

# Row 49
Set-TextValue "E49" "  +0.18%  "

# Row 50
Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0958"
Set-TextValue "E50" "  -2.52%  "

# Row 51
Set-TextValue "B51" "Mantle"
Set-TextValue "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D51" "0.406"
Set-TextValue "E51" "  -0.65%  "
